$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (BF) was populated from the source workbook's file name
# ("5-2-2013-14"), which is wrong because of how the NBA stats site lays out
# its dates. Replace it with the correct ISO date "2014-05-02" on every data
# row (BF2:BF31). A leading apostrophe forces Excel to keep the value as
# literal text instead of re-interpreting the ISO-looking string as a date
# serial number.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)
    if ($cell.Value() -eq "5-2-2013-14") {
        $cell.Value = "'2014-05-02"
    }
}
